# actualizacion a Agosto 2022
# Adds the August 2022 (period 8) data row to Sheet1, row 57, mirroring
# the layout/formatting of the existing rows, and updates the sheet's
# active selection / scroll position to reflect where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$row = 57

# --- Plain (unformatted / General) columns -------------------------------
$ws.Cells.Item($row, 1).Value  = 2022          # A57
$ws.Cells.Item($row, 2).Value  = 8              # B57
$ws.Cells.Item($row, 3).Value  = 76421          # C57
$ws.Cells.Item($row, 15).Value = 0.18907420735138236   # O57
$ws.Cells.Item($row, 17).Value = 277166.14269301592    # Q57
$ws.Cells.Item($row, 18).Value = 0.55433228538603185   # R57

# --- Columns formatted with the "0.00" numeric style (style index 1) -----
$fmtCols = @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16)  # D..N, P
$fmtVals = @{
    4  = 126858.86
    5  = 76.5
    6  = 127005.36
    7  = 49170.55
    8  = 44880.77
    9  = 7139
    10 = 2385.8000000000002
    11 = 8980
    12 = 112556.12
    13 = 14449.24
    14 = 2076.04
    16 = 500000
}

foreach ($col in $fmtCols) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $fmtVals[$col]
    $cell.NumberFormat = "0.00"
}

# --- Update the view: window scrolls so that L43 becomes the top-left ---
# --- visible cell, and the active cell / selection moves to N56 ---------
$excel.Goto($ws.Range("L43"), $true)
$ws.Range("N56").Select()
